$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 0. Helper functions
# ---------------------------------------------------------------------
# Copy the (border/shading/alignment) style from $src onto $dst, then
# overwrite $dst's value - used for the bordered "s=2" header / index
# cells so the new cells carry the same look as their neighbours.
function Set-StyledValue($src, $dst, $value) {
    $src.Copy($dst)
    $dst.Value = $value
}

# Write $value into $range as literal TEXT (not a number), mirroring
# source cells such as "7.16" that are stored as inline strings even
# though they look numeric. NumberFormat forces Excel to keep the
# literal text instead of parsing it as a number; ClearFormats()
# afterwards drops that temporary formatting again so the cell is left
# with the default (unstyled) look, matching the source file.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q3" sheet right after "总计" (i.e. before the
#    current second sheet, "2021-Q2"). Duplicating the "2021-Q2" sheet
#    (rather than Worksheets.Add()) keeps the sheetPr/outline settings
#    that sheet carries, then we wipe it completely and rebuild it.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)
$q2Sheet    = $wb.Worksheets.Item(2)
$q2Sheet.Copy($q2Sheet)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q3"
$newSheet.Cells.Clear()

# ---------------------------------------------------------------------
# 2. Populate the new "2022-Q3" sheet with the fund-holding data.
# ---------------------------------------------------------------------
Set-StyledValue $totalSheet.Range("B1") $newSheet.Range("B1") "基金代码"
Set-StyledValue $totalSheet.Range("B1") $newSheet.Range("C1") "基金名称"
Set-StyledValue $totalSheet.Range("B1") $newSheet.Range("D1") "基金规模"
Set-StyledValue $totalSheet.Range("B1") $newSheet.Range("E1") "股票总仓位"
Set-StyledValue $totalSheet.Range("B1") $newSheet.Range("F1") "仓位占比"
Set-StyledValue $totalSheet.Range("B1") $newSheet.Range("G1") "持有市值(亿元)"
Set-StyledValue $totalSheet.Range("B1") $newSheet.Range("H1") "仓位排名"

Set-StyledValue $totalSheet.Range("A2") $newSheet.Range("A2") 0
Set-TextValue $newSheet.Range("B2") "970042"
Set-TextValue $newSheet.Range("C2") "国海量化优选一年持有股票C"
Set-TextValue $newSheet.Range("D2") "7.16"
Set-TextValue $newSheet.Range("E2") "87.31"
Set-TextValue $newSheet.Range("F2") "0.34"
Set-TextValue $newSheet.Range("G2") "0.0243"
$newSheet.Range("H2").Value = 7

Set-StyledValue $totalSheet.Range("A2") $newSheet.Range("A3") 1
Set-TextValue $newSheet.Range("B3") "970041"
Set-TextValue $newSheet.Range("C3") "国海量化优选一年持有股票A"
Set-TextValue $newSheet.Range("D3") "0.63"
Set-TextValue $newSheet.Range("E3") "87.31"
Set-TextValue $newSheet.Range("F3") "0.34"
Set-TextValue $newSheet.Range("G3") "0.0021"
$newSheet.Range("H3").Value = 7

# ---------------------------------------------------------------------
# 3. Update the "总计" summary sheet: shift the existing two rows down
#    and insert the new 2022-Q3 figures on top.
# ---------------------------------------------------------------------
Set-StyledValue $totalSheet.Range("A3") $totalSheet.Range("A4") 2
$totalSheet.Range("B4").Value = "2021-Q1"
$totalSheet.Range("C4").Value = 3
$totalSheet.Range("D4").Value = 0.06

$totalSheet.Range("B3").Value = "2021-Q2"
$totalSheet.Range("C3").Value = 3
$totalSheet.Range("D3").Value = 0.09

$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.03

# ---------------------------------------------------------------------
# 4. Keep the original active/selected tab on "2021-Q1" (it was the
#    selected sheet before the edit; inserting a sheet would otherwise
#    steal the selection).
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q1").Activate()
